# Generate Report for Archive
#
# The "bd241674-a7b5-470e-bba4-68335070b746" file's status moved from
# "Ready for handoff" to "In Translation", and as the report was
# regenerated its row now sorts ahead of the "13855b4a-e3d5-4cc6-8ba6-
# decf47c1a37c" row (rows 5 and 6 swap their data) on every sheet. The
# hyperlink targets (Address / r:id) stay pinned to their original
# relationship - only the cell values and the hyperlinks' displayed text
# change.

$wb = $excel.ActiveWorkbook

function Set-CellAndLink {
    param(
        $ws,
        [string]$cellAddr,
        [string]$newValue,
        [bool]$isHyperlink
    )

    $target = $ws.Range($cellAddr)
    $target.Value = $newValue

    if ($isHyperlink) {
        $targetAddr = $target.Address()
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Address() -eq $targetAddr) {
                $h.TextToDisplay = $newValue
            }
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndLink $wsOverview "A5" "bd241674-a7b5-470e-bba4-68335070b746.md" $true
$wsOverview.Range("B5").Value = "In Translation"
$wsOverview.Range("C5").Value = "In Translation"
$wsOverview.Range("D5").Value = "2016-03-21 18:39:10"

Set-CellAndLink $wsOverview "A6" "13855b4a-e3d5-4cc6-8ba6-decf47c1a37c.md" $true
$wsOverview.Range("D6").Value = "2016-03-21 18:36:43"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndLink $wsZhCn "A5" "bd241674-a7b5-470e-bba4-68335070b746.md" $true
$wsZhCn.Range("C5").Value = "In Translation"
Set-CellAndLink $wsZhCn "D5" "bd241674-a7b5-470e-bba4-68335070b746.58424105711c9eb457a1a03d7153a9eacd355735.zh-cn.xlf" $true
$wsZhCn.Range("E5").Value = "2016-03-21 18:39:07"

Set-CellAndLink $wsZhCn "A6" "13855b4a-e3d5-4cc6-8ba6-decf47c1a37c.md" $true
Set-CellAndLink $wsZhCn "D6" "13855b4a-e3d5-4cc6-8ba6-decf47c1a37c.d0b1200a2aed74349ffae9c6da5c6d7f7e7e477d.zh-cn.xlf" $true
$wsZhCn.Range("E6").Value = "2016-03-21 18:36:34"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndLink $wsDeDe "A5" "bd241674-a7b5-470e-bba4-68335070b746.md" $true
$wsDeDe.Range("C5").Value = "In Translation"
Set-CellAndLink $wsDeDe "D5" "bd241674-a7b5-470e-bba4-68335070b746.58424105711c9eb457a1a03d7153a9eacd355735.de-de.xlf" $true
$wsDeDe.Range("E5").Value = "2016-03-21 18:39:10"

Set-CellAndLink $wsDeDe "A6" "13855b4a-e3d5-4cc6-8ba6-decf47c1a37c.md" $true
Set-CellAndLink $wsDeDe "D6" "13855b4a-e3d5-4cc6-8ba6-decf47c1a37c.d0b1200a2aed74349ffae9c6da5c6d7f7e7e477d.de-de.xlf" $true
$wsDeDe.Range("E6").Value = "2016-03-21 18:36:43"
